$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 38 - this shifts existing rows 38:144 down to 39:145
$ws.Rows("38:38").Insert()

# Populate the newly inserted row 38 with the new weekly data point
$ws.Range("A38").Value = 5
$ws.Range("B38").Value = "Macroferia Regional de Talca"
$ws.Range("C38").Value = "Maule"
$ws.Range("D38").Value = 44414
$ws.Range("E38").Value = 7
$ws.Range("F38").Value = 100114014
$ws.Range("G38").Value = "Betarraga"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 3000
$ws.Range("K38").Value = 600
$ws.Range("L38").Value = 600
$ws.Range("M38").Value = 600
$ws.Range("N38").Value = "$/paquete 5 unidades"
$ws.Range("O38").Value = "Región del Maule"
$ws.Range("P38").Value = 120
$ws.Range("Q38").Value = 5
$ws.Range("R38").Value = "Hortaliza"
